$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.117.74"
$ws.Range("E2").Value = "  -2.63%  "
$ws.Range("D3").Value = "1.867.80"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'306.91"
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.5108"
$ws.Range("E7").Value = "  +2.14%  "
$ws.Range("D8").Value = "'0.3747"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").Value = "'0.07125"
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("D10").Value = "'0.8869"
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("D11").Value = "'20.59"
$ws.Range("E11").Value = "  -2.91%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.864.04"
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07538"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").Value = "'5.323"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").Value = "'89.03"
$ws.Range("E15").Value = "  -3.47%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "'0.000008464"
$ws.Range("E17").Value = "  -2.92%  "
$ws.Range("D18").Value = "'14.10"
$ws.Range("E18").Value = "  -3.52%  "
$ws.Range("D19").Value = "'0.9996"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "27.171.39"
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("D21").Value = "'5.052"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("D22").Value = "2.108.73"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").Value = "'10.55"
$ws.Range("E23").Value = "  -2.79%  "
$ws.Range("D24").Value = "'6.462"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").Value = "'149.75"
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("D26").Value = "'1.841"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D28").Value = "'2.096"
$ws.Range("E28").Value = "  -5.44%  "
$ws.Range("D29").Value = "'113.04"
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("D30").Value = "'4.706"
$ws.Range("E30").Value = "  -3.74%  "
$ws.Range("D31").Value = "'4.663"
$ws.Range("E31").Value = "  -3.14%  "
$ws.Range("D32").Value = "'0.09022"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "'0.05135"
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("D34").Value = "'3.095"
$ws.Range("E34").Value = "  -3.37%  "
$ws.Range("D35").Value = "'1.153"
$ws.Range("E35").Value = "  -6.37%  "
$ws.Range("D36").Value = "'0.7339"
$ws.Range("E36").Value = "  -6.42%  "
$ws.Range("D37").Value = "'0.02057"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").Value = "'2.511"
$ws.Range("E38").Value = "  -5.08%  "
$ws.Range("D39").Value = "'3.061"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "'1.072"
$ws.Range("E40").Value = "  -1.97%  "
$ws.Range("D41").Value = "'0.5337"
$ws.Range("E41").Value = "  -3.28%  "
$ws.Range("D42").Value = "'6.568"
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("D43").Value = "'116.40"
$ws.Range("E43").Value = "  +1.75%  "
$ws.Range("D44").Value = "'8.330"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("D45").Value = "'0.1470"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.4622"
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").Value = "'10.01"
$ws.Range("E48").Value = "  -4.82%  "
$ws.Range("D49").Value = "'1.568"
$ws.Range("E49").Value = "  -4.04%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'36.63"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'64.25"
$ws.Range("E51").Value = "  -4.34%  "
